$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Intitulé du champ" column (C) for the existing "utilisateur" rows ---
$ws.Range("C3").Value = "identifiant de l'utilisateur"
$ws.Range("C4").Value = "pseudo de l'utilisateur"
$ws.Range("C5").Value = "mot de passe de l'utilisateur"
$ws.Range("C6").Value = "date inscription de l'utilisateur"
$ws.Range("C7").Value = "nom de l'utilisateur"
$ws.Range("C8").Value = "prénom  de l'utilisateur"
$ws.Range("C9").Value = "age de l'utilisateur"
$ws.Range("C10").Value = "sexe de l'utilisateur"
$ws.Range("C11").Value = "email  de l'utilisateur"

# --- Row 12 changes: tel -> tel_fixe, plus its new "Intitulé" ---
$ws.Range("A12").Value = "tel_fixe"
$ws.Range("C12").Value = "téléhpone de l'utilisateur"

# --- Insert 4 new rows (13-16) for tel_mobile / adresse / cp / ville ---
$ws.Rows("13:16").Insert()

$ws.Range("A13").Value = "tel_mobile"
$ws.Range("B13").Value = "VARCHAR(30)"

$ws.Range("A14").Value = "adresse"
$ws.Range("B14").Value = "VARCHAR(120)"

$ws.Range("A15").Value = "cp"
$ws.Range("B15").Value = "INT"

$ws.Range("A16").Value = "ville"
$ws.Range("B16").Value = "VARCHAR(50)"

# --- Column C width change ---
$ws.Columns("C").ColumnWidth = 35.8

# --- Sheet view: scroll position and selection ---
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("K25").Select()
